$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''28.524.28'
$ws.Range('E2').Value = '  +1.08%  '
$ws.Range('D3').Value = '''1.912.48'
$ws.Range('E3').Value = '  +4.40%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''315.19'
$ws.Range('E5').Value = '  +1.57%  '
$ws.Range('D6').Value = '''1.000'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = '''0.5144'
$ws.Range('E7').Value = '  +3.43%  '
$ws.Range('D8').Value = '''0.3974'
$ws.Range('E8').Value = '  +1.20%  '
$ws.Range('D9').Value = '''0.09812'
$ws.Range('E9').Value = '  -2.79%  '
$ws.Range('E10').Value = '  +3.35%  '
$ws.Range('D11').Value = '''42.25'
$ws.Range('E11').Value = '  +2.53%  '
$ws.Range('D12').Value = '''6.537'
$ws.Range('E12').Value = '  +1.50%  '
$ws.Range('D13').Value = '''21.19'
$ws.Range('E13').Value = '  +2.34%  '
$ws.Range('D14').Value = '''1.911.37'
$ws.Range('E14').Value = '  +4.85%  '
$ws.Range('D15').Value = '''7.488'
$ws.Range('E15').Value = '  +1.98%  '
$ws.Range('D16').Value = '''1.001'
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').Value = '''94.61'
$ws.Range('E18').Value = '  -1.04%  '
$ws.Range('D19').Value = '''0.06661'
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').Value = '''18.26'
$ws.Range('E20').Value = '  +5.73%  '
$ws.Range('D21').Value = '''1.000'
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('D22').Value = '''6.322'
$ws.Range('E22').Value = '  +5.05%  '
$ws.Range('D23').Value = '''28.579.93'
$ws.Range('E23').Value = '  +1.10%  '
$ws.Range('D24').Value = '''11.50'
$ws.Range('E24').Value = '  +1.54%  '
$ws.Range('D25').Value = '''2.320'
$ws.Range('E25').Value = '  +4.03%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').Value = '''2.678'
$ws.Range('E26').Value = '  +9.92%  '
$ws.Range('B27').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C27').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D27').Value = '''2.127.17'
$ws.Range('E27').Value = '  +4.24%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '''21.27'
$ws.Range('E28').Value = '  +2.16%  '
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').Value = '''157.55'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '''128.97'
$ws.Range('E30').Value = '  +1.66%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '''1.114'
$ws.Range('E31').Value = '  +7.18%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '''0.1079'
$ws.Range('E32').Value = '  +2.60%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '''5.759'
$ws.Range('E33').Value = '  +2.70%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '''3.629'
$ws.Range('E34').Value = '  +0.83%  '
$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D35').Value = '''9.880'
$ws.Range('E35').Value = '  +8.95%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '''0.06790'
$ws.Range('E36').Value = '  +0.31%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '''0.02434'
$ws.Range('E37').Value = '  +3.28%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').Value = '''1.274'
$ws.Range('E38').Value = '  +7.11%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').Value = '''0.2213'
$ws.Range('E39').Value = '  +2.71%  '
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').Value = '''11.81'
$ws.Range('E40').Value = '  +2.96%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '''0.6483'
$ws.Range('E41').Value = '  +3.91%  '
$ws.Range('B42').Value = 'InternetComputer(DFINITY)'
$ws.Range('C42').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D42').Value = '''5.091'
$ws.Range('E42').Value = '  +2.02%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '''1.188'
$ws.Range('E43').Value = '  +1.13%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').Value = '''1.000'
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '''13.62'
$ws.Range('E45').Value = '  +2.57%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '''0.6105'
$ws.Range('E46').Value = '  +2.65%  '
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').Value = '''3.778'
$ws.Range('E47').Value = '  +2.42%  '
$ws.Range('B48').Value = 'WEMIXTOKEN'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '''1.293'
$ws.Range('E48').Value = '  +1.63%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '''2.045'
$ws.Range('E49').Value = '  +4.65%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = '''124.94'
$ws.Range('E50').Value = '  +0.64%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').Value = '''1.204'
$ws.Range('E51').Value = '  +1.85%  '
